$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.254529666666667
$ws.Range("H2").Value = 3.763589
$ws.Range("I2").Value = 0.01270475613604028
$ws.Range("J2").Value = 0.01270475613604028
$ws.Range("M2").Value = 35.04689966666667
$ws.Range("N2").Value = 105.140699
$ws.Range("O2").Value = 0.3824629895491901
$ws.Range("P2").Value = 0.3824629895491901
$ws.Range("Q2").Value = 43.96737535652345
$ws.Range("R2").Value = 395.706378208711
$ws.Range("S2").Value = 0.004859099013283382
$ws.Range("T2").Value = 0.00485909901328338
$ws.Range("G3").Value = 1.254529666666667
$ws.Range("H3").Value = 3.763589
$ws.Range("I3").Value = 0.01270475613604028
$ws.Range("J3").Value = 0.01270475613604028
$ws.Range("O3").Value = 0.3264402385872224
$ws.Range("P3").Value = 0.3264402385872223
$ws.Range("Q3").Value = 37.52708338748034
$ws.Range("R3").Value = 337.743750487323
$ws.Range("S3").Value = 0.004147343624241465
$ws.Range("T3").Value = 0.004147343624241464
$ws.Range("G4").Value = 1.254529666666667
$ws.Range("H4").Value = 3.763589
$ws.Range("I4").Value = 0.01270475613604028
$ws.Range("J4").Value = 0.01270475613604028
$ws.Range("M4").Value = 8.911727666666666
$ws.Range("N4").Value = 26.735183
$ws.Range("O4").Value = 0.09725271102035077
$ws.Range("P4").Value = 0.09725271102035075
$ws.Range("Q4").Value = 11.18002673908744
$ws.Range("R4").Value = 100.620240651787
$ws.Range("S4").Value = 0.001235571977082353
$ws.Range("T4").Value = 0.001235571977082353
$ws.Range("G5").Value = 1.254529666666667
$ws.Range("H5").Value = 3.763589
$ws.Range("I5").Value = 0.01270475613604028
$ws.Range("J5").Value = 0.01270475613604028
$ws.Range("M5").Value = 17.76285166666667
$ws.Range("N5").Value = 53.288555
$ws.Range("O5").Value = 0.1938440608432367
$ws.Range("P5").Value = 0.1938440608432367
$ws.Range("Q5").Value = 22.28402438043278
$ws.Range("R5").Value = 200.556219423895
$ws.Range("S5").Value = 0.002462741521433077
$ws.Range("T5").Value = 0.002462741521433076
$ws.Range("I6").Value = 0.9734152842234517
$ws.Range("J6").Value = 0.9734152842234516
$ws.Range("M6").Value = 35.04689966666667
$ws.Range("N6").Value = 105.140699
$ws.Range("O6").Value = 0.3824629895491901
$ws.Range("P6").Value = 0.3824629895491901
$ws.Range("Q6").Value = 3368.700250595174
$ws.Range("R6").Value = 30318.30225535657
$ws.Range("S6").Value = 0.3722953196769759
$ws.Range("T6").Value = 0.3722953196769758
$ws.Range("I7").Value = 0.9734152842234517
$ws.Range("J7").Value = 0.9734152842234516
$ws.Range("O7").Value = 0.3264402385872224
$ws.Range("P7").Value = 0.3264402385872223
$ws.Range("S7").Value = 0.3177619176263525
$ws.Range("T7").Value = 0.3177619176263524
$ws.Range("I8").Value = 0.9734152842234517
$ws.Range("J8").Value = 0.9734152842234516
$ws.Range("M8").Value = 8.911727666666666
$ws.Range("N8").Value = 26.735183
$ws.Range("O8").Value = 0.09725271102035077
$ws.Range("P8").Value = 0.09725271102035075
$ws.Range("Q8").Value = 856.5932938281857
$ws.Range("R8").Value = 7709.339644453672
$ws.Range("S8").Value = 0.09466727533937595
$ws.Range("T8").Value = 0.09466727533937593
$ws.Range("I9").Value = 0.9734152842234517
$ws.Range("J9").Value = 0.9734152842234516
$ws.Range("M9").Value = 17.76285166666667
$ws.Range("N9").Value = 53.288555
$ws.Range("O9").Value = 0.1938440608432367
$ws.Range("P9").Value = 0.1938440608432367
$ws.Range("Q9").Value = 1707.361376609782
$ws.Range("R9").Value = 15366.25238948804
$ws.Range("S9").Value = 0.1886907715807473
$ws.Range("T9").Value = 0.1886907715807473
$ws.Range("G10").Value = 1.151276666666667
$ws.Range("H10").Value = 3.45383
$ws.Range("I10").Value = 0.01165910195968263
$ws.Range("J10").Value = 0.01165910195968263
$ws.Range("M10").Value = 35.04689966666667
$ws.Range("N10").Value = 105.140699
$ws.Range("O10").Value = 0.3824629895491901
$ws.Range("P10").Value = 0.3824629895491901
$ws.Range("Q10").Value = 40.34867782524112
$ws.Range("R10").Value = 363.1381004271701
$ws.Range("S10").Value = 0.004459174990959039
$ws.Range("T10").Value = 0.004459174990959039
$ws.Range("G11").Value = 1.151276666666667
$ws.Range("H11").Value = 3.45383
$ws.Range("I11").Value = 0.01165910195968263
$ws.Range("J11").Value = 0.01165910195968263
$ws.Range("O11").Value = 0.3264402385872224
$ws.Range("P11").Value = 0.3264402385872223
$ws.Range("Q11").Value = 34.43844862342333
$ws.Range("R11").Value = 309.94603761081
$ws.Range("S11").Value = 0.00380600002543155
$ws.Range("T11").Value = 0.003806000025431549
$ws.Range("G12").Value = 1.151276666666667
$ws.Range("H12").Value = 3.45383
$ws.Range("I12").Value = 0.01165910195968263
$ws.Range("J12").Value = 0.01165910195968263
$ws.Range("M12").Value = 8.911727666666666
$ws.Range("N12").Value = 26.735183
$ws.Range("O12").Value = 0.09725271102035077
$ws.Range("P12").Value = 0.09725271102035075
$ws.Range("Q12").Value = 10.25986412232111
$ws.Range("R12").Value = 92.33877710089
$ws.Range("S12").Value = 0.00113387927364182
$ws.Range("T12").Value = 0.00113387927364182
$ws.Range("G13").Value = 1.151276666666667
$ws.Range("H13").Value = 3.45383
$ws.Range("I13").Value = 0.01165910195968263
$ws.Range("J13").Value = 0.01165910195968263
$ws.Range("M13").Value = 17.76285166666667
$ws.Range("N13").Value = 53.288555
$ws.Range("O13").Value = 0.1938440608432367
$ws.Range("P13").Value = 0.1938440608432367
$ws.Range("Q13").Value = 20.44995665729445
$ws.Range("R13").Value = 184.04960991565
$ws.Range("S13").Value = 0.002260047669650221
$ws.Range("T13").Value = 0.00226004766965022
$ws.Range("G14").Value = 0.2192983333333333
$ws.Range("H14").Value = 0.6578949999999999
$ws.Range("I14").Value = 0.002220857680825461
$ws.Range("J14").Value = 0.002220857680825461
$ws.Range("M14").Value = 35.04689966666667
$ws.Range("N14").Value = 105.140699
$ws.Range("O14").Value = 0.3824629895491901
$ws.Range("P14").Value = 0.3824629895491901
$ws.Range("Q14").Value = 7.685726685400555
$ws.Range("R14").Value = 69.17154016860499
$ws.Range("S14").Value = 0.000849395867971787
$ws.Range("T14").Value = 0.0008493958679717866
$ws.Range("G15").Value = 0.2192983333333333
$ws.Range("H15").Value = 0.6578949999999999
$ws.Range("I15").Value = 0.002220857680825461
$ws.Range("J15").Value = 0.002220857680825461
$ws.Range("O15").Value = 0.3264402385872224
$ws.Range("P15").Value = 0.3264402385872223
$ws.Range("Q15").Value = 6.559930036251665
$ws.Range("R15").Value = 59.03937032626499
$ws.Range("S15").Value = 0.0007249773111969289
$ws.Range("T15").Value = 0.0007249773111969287
$ws.Range("G16").Value = 0.2192983333333333
$ws.Range("H16").Value = 0.6578949999999999
$ws.Range("I16").Value = 0.002220857680825461
$ws.Range("J16").Value = 0.002220857680825461
$ws.Range("M16").Value = 8.911727666666666
$ws.Range("N16").Value = 26.735183
$ws.Range("O16").Value = 0.09725271102035077
$ws.Range("P16").Value = 0.09725271102035075
$ws.Range("Q16").Value = 1.954327024420555
$ws.Range("R16").Value = 17.588943219785
$ws.Range("S16").Value = 0.000215984430250645
$ws.Range("T16").Value = 0.0002159844302506449
$ws.Range("G17").Value = 0.2192983333333333
$ws.Range("H17").Value = 0.6578949999999999
$ws.Range("I17").Value = 0.002220857680825461
$ws.Range("J17").Value = 0.002220857680825461
$ws.Range("M17").Value = 17.76285166666667
$ws.Range("N17").Value = 53.288555
$ws.Range("O17").Value = 0.1938440608432367
$ws.Range("P17").Value = 0.1938440608432367
$ws.Range("Q17").Value = 3.895363765747221
$ws.Range("R17").Value = 35.05827389172499
$ws.Range("S17").Value = 0.0004305000714061004
$ws.Range("T17").Value = 0.0004305000714061002
